# Regenerate save_data: replace column G ("K", strikeout count) values
# with newly computed figures for rows 2-48 (row 1 is the header row).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of row -> new K value, taken from the recomputed save_data.
$kValues = [ordered]@{
    2  = 2
    3  = 3
    4  = 0
    5  = 1
    6  = 1
    7  = 0
    8  = 3
    9  = 2
    10 = 1
    11 = 0
    12 = 1
    13 = 1
    14 = 0
    15 = 3
    16 = 2
    17 = 0
    18 = 1
    19 = 2
    20 = 1
    21 = 1
    22 = 0
    23 = 1
    24 = 0
    25 = 2
    26 = 2
    27 = 1
    28 = 0
    29 = 0
    30 = 1
    31 = 0
    32 = 0
    33 = 1
    34 = 2
    35 = 1
    37 = 1
    38 = 1
    39 = 1
    40 = 0
    41 = 1
    42 = 2
    43 = 1
    44 = 0
    45 = 2
    46 = 1
    47 = 0
    48 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
